# Update the date heading.
$d = $word.ActiveDocument
$d.Content.Find.Execute("2025-09-07 Sunday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2025-09-08 Monday", 2)

# Update each practice-problem cell in the table, addressed directly by
# (row, column) so that duplicate/overlapping old-vs-new values (e.g. the
# "86÷6=14, 2" that is simultaneously a target value in one cell and a
# source value in another) cannot cross-contaminate each other the way a
# blanket Find/Replace across the whole document could.
$tbl = $d.Tables.Item(1)

$tbl.Cell(1, 1).Range.Text = "46÷4=11, 2"
$tbl.Cell(1, 2).Range.Text = "87÷6=14, 3"
$tbl.Cell(1, 3).Range.Text = "57÷7=8, 1"
$tbl.Cell(1, 4).Range.Text = "14÷6=2, 2"
$tbl.Cell(1, 5).Range.Text = "82÷6=13, 4"

$tbl.Cell(5, 1).Range.Text = "35÷3=11, 2"
$tbl.Cell(5, 2).Range.Text = "33÷6=5, 3"
$tbl.Cell(5, 3).Range.Text = "14÷3=4, 2"
$tbl.Cell(5, 4).Range.Text = "83÷9=9, 2"
$tbl.Cell(5, 5).Range.Text = "73÷8=9, 1"

$tbl.Cell(9, 1).Range.Text = "43÷7=6, 1"
$tbl.Cell(9, 2).Range.Text = "27÷5=5, 2"
$tbl.Cell(9, 3).Range.Text = "83÷8=10, 3"
$tbl.Cell(9, 4).Range.Text = "50÷3=16, 2"
$tbl.Cell(9, 5).Range.Text = "51÷5=10, 1"

$tbl.Cell(13, 1).Range.Text = "55÷3=18, 1"
$tbl.Cell(13, 2).Range.Text = "69÷9=7, 6"
$tbl.Cell(13, 3).Range.Text = "86÷6=14, 2"
$tbl.Cell(13, 4).Range.Text = "46÷6=7, 4"
$tbl.Cell(13, 5).Range.Text = "43÷7=6, 1"

$tbl.Cell(17, 1).Range.Text = "64÷9=7, 1"
$tbl.Cell(17, 2).Range.Text = "70÷5=14, 0"
$tbl.Cell(17, 3).Range.Text = "36÷9=4, 0"
$tbl.Cell(17, 4).Range.Text = "59÷4=14, 3"
$tbl.Cell(17, 5).Range.Text = "67÷6=11, 1"
